$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a brand-new "2022-Q3" sheet right after "总计" (position 2),
#    pushing every existing quarter sheet one slot to the right.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

# Helper style source: a cell that already carries the shared "bold + boxed
# header / index" style (style index 2 in the original workbook) so the new
# sheet visually matches its siblings.
$styleSrc = $total.Range("B1")
$idxStyleSrc = $total.Range("A2")

# ---- header row ------------------------------------------------------------
$styleSrc.Copy()
$q3.Range("B1:D1").PasteSpecial(-4122)
$q3.Range("E1:H1").PasteSpecial(-4122)

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# ---- index column (A2:A13) -------------------------------------------------
$idxStyleSrc.Copy()
$q3.Range("A2:A13").PasteSpecial(-4122)

# ---- data rows --------------------------------------------------------------
# columns B,D,E,F,G must stay text (leading zeros / trailing zeros matter),
# column H is a genuine number.
$textCols = @("B", "D", "E", "F", "G")
foreach ($col in $textCols) {
    $q3.Range("$col 2:$col 13".Replace(" ", "")).NumberFormat = "@"
}

$rows = @(
    @(0, "410003", "华富成长趋势混合", "10.39", "93.26", "4.81", "0.4998", 7),
    @(1, "410007", "华富价值增长混合", "8.46", "79.64", "3.77", "0.3189", 9),
    @(2, "001437", "易方达瑞享灵活配置混合I", "4.27", "92.06", "5.06", "0.2161", 10),
    @(3, "001438", "易方达瑞享灵活配置混合E", "4.27", "92.06", "5.06", "0.2161", 10),
    @(4, "009398", "华富成长企业精选股票", "4.15", "92.29", "3.78", "0.1569", 8),
    @(5, "002581", "招商丰凯灵活配置混合A", "4.22", "48.56", "1.72", "0.0726", 5),
    @(6, "012586", "南方港股创新视野一年持有混合A", "2.16", "50.74", "3.11", "0.0672", 7),
    @(7, "519644", "银河智联主题灵活配置混合", "0.90", "93.05", "4.95", "0.0446", 6),
    @(8, "003152", "华富天鑫灵活配置混合A", "0.88", "94.01", "4.33", "0.0381", 8),
    @(9, "002582", "招商丰凯灵活配置混合C", "1.15", "48.56", "1.72", "0.0198", 5),
    @(10, "012587", "南方港股创新视野一年持有混合C", "0.19", "50.74", "3.11", "0.0059", 7),
    @(11, "003153", "华富天鑫灵活配置混合C", "0.06", "94.01", "4.33", "0.0026", 8)
)

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2) Update "总计" (sheet1): shift the quarterly summary rows down by one and
#    fill in the brand-new 2022-Q3 figures. The index column (A) is a plain
#    0-based row counter, so it is left untouched for existing rows; only
#    the new last row (row 9) needs a fresh index value.
# ---------------------------------------------------------------------------
$total.Range("A8").Copy()
$total.Range("A9").PasteSpecial(-4122)
$total.Range("A9").Value = 7

$summary = @(
    @("2022-Q3", 12, 1.66),
    @("2022-Q2", 23, 3.97),
    @("2022-Q1", 15, 6.23),
    @("2021-Q4", 21, 14.11),
    @("2021-Q3", 13, 10.58),
    @("2021-Q2", 12, 13.91),
    @("2021-Q1", 9, 6.6),
    @("2020-Q4", 8, 4.52)
)

$r = 2
foreach ($row in $summary) {
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}
